$wb = $excel.ActiveWorkbook

$wsDevice = $wb.Worksheets.Item("DEVICE_CONFIG")
$wsNet = $wb.Worksheets.Item("BASIC_NETWORK")

# --- BASIC_NETWORK sheet: finalize values for a single device/interface ---
$wsNet.Range("J2").Value = "true, false"

# --- DEVICE_CONFIG sheet: device name filled in with the real (single) device value ---
$wsDevice.Range("D2").Value = "new22"

$wsNet.Range("L2").Value = "192.168.1.1"

# Row 3 becomes a blank spacer row spanning every used column (A:N). The
# cells that weren't part of the row yet (B3, C3, F3:N3) pick up the same
# formatting already used by the rest of the row (A3/D3/E3).
$wsNet.Range("A3").Copy()
$wsNet.Range("B3:C3").PasteSpecial(-4122)
$wsNet.Range("F3:N3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the now-unneeded "true/false" dropdown that used to live on D3 (the
# template row), and shrink the enabled/disabled dropdown so it no longer
# covers row 3 now that row 3 is just a blank spacer.
$wsNet.Range("D3").Validation.Delete()
$wsNet.Range("E3:E17").Validation.Delete()
$wsNet.Range("E4:E17").Validation.Add(3, 1, 1, '"enabled, disabled"')

# Selection / active-sheet bookkeeping: BASIC_NETWORK is no longer the active
# tab, selection moves to L2; DEVICE_CONFIG becomes the active tab, selection
# moves to D2.
$wsNet.Range("L2").Select()

$wsDevice.Activate()
$wsDevice.Range("D2").Select()
